$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (2020-10-03) appended below the existing table, matching
# columns A:DX (date serial + 127 numeric columns).
$rowValues = @(
    44107,848147,2746,117695,67859,273481,29376,7589,6360,8730,9911,20720,4030,
    24291,34097,8280,12249,15533,15631,18732,16305,3863,4113,11441,32249,14158,
    12644,64462,2581,1505,836,480,950,556,857,2087,6142,38335,10481,2580,49533,
    1175,23442,1539,10731,1688,1618,8825,2075,966,2504,2705,67987,14372,7379,
    10160,7657,257,1471,2762,746,2187,10103,9623,10863,14428,1977,913,14558,
    11898,13987,3559,2416,6513,5204,2889,6358,4033,2418,1339,3106,2273,2206,
    2102,6766,2335,1552,1914,2202,2370,2827,1906,1245,1242,1198,3511,1552,1007,
    1227,1785,1660,824,947,1422,1888,1771,1744,1339,337,373,841,796,518,544,
    394,702,767,534,515,374,527,141768,357325,22187,156134,97135,48603,13543
)

$targetRow = 213
$arr = New-Object 'object[,]' 1,$rowValues.Length
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $arr[0,$i] = $rowValues[$i]
}

$startCell = $ws.Cells.Item($targetRow, 1)
$endCell = $ws.Cells.Item($targetRow, $rowValues.Length)
$ws.Range($startCell, $endCell).Value = $arr

# Scroll the frozen (bottom-right) pane so its top-left visible cell is
# B181, then move the active selection down to A209, matching where the
# author left the view after entering the new data.
$win = $excel.ActiveWindow
$win.ScrollRow = 181
$win.ScrollColumn = 2
$ws.Range("A209").Select()
